$wb = $excel.ActiveWorkbook

# Sheet "zh-cn": update Correspond Handoff Datetime (E2) and Correspond Handback DateTime (H2)
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E2").Value = "2016-03-12 00:45:54"
$wsZh.Range("H2").Value = "2016-03-12 00:46:20"

# Sheet "de-de": update Correspond Handoff Datetime (E2) and Correspond Handback DateTime (H2)
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E2").Value = "2016-03-12 00:45:58"
$wsDe.Range("H2").Value = "2016-03-12 00:46:25"
